$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(19).Delete()
$ws.Range("B25").Select()
